$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Stop time (E12): 17:30 -> 17:20
$ws.Range("E12").Value = (17*60 + 20) / 1440.0

# Update Delta Time minutes (G12): 190 -> 180
$ws.Range("G12").Value = 180

# Update selected/active cell to G16
$ws.Range("G16").Select()
